$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 5.515135932605308

$ws.Range("C3").Value = -2.764179109379705
$ws.Range("E3").Value = -0.2596132895792413

$ws.Range("C4").Value = 3.38265053313096
$ws.Range("E4").Value = 1.421687849828523

$ws.Range("C5").Value = 5.813045170083808
$ws.Range("E5").Value = 2.371160938652705

$ws.Range("C6").Value = 4.597183386292891
$ws.Range("E6").Value = 6.213583554874536

$ws.Range("C7").Value = 2.113053977048707
$ws.Range("E7").Value = 3.80132608443593

$ws.Range("C8").Value = 5.262295419893648
$ws.Range("E8").Value = 3.911207409579687

$ws.Range("C9").Value = 4.78031692483154
$ws.Range("E9").Value = 4.481753591536197

$ws.Range("C10").Value = 5.386945580119185
$ws.Range("E10").Value = 5.291601650505706

$ws.Range("C11").Value = 4.884223728030879
$ws.Range("E11").Value = 4.848884192354119

$ws.Range("C12").Value = 6.072005530313129
$ws.Range("E12").Value = 5.47980442665561

$ws.Range("C13").Value = 4.927928448556984
$ws.Range("E13").Value = 5.802765867180804

$ws.Range("C14").Value = 2.502519143054571
$ws.Range("E14").Value = 3.379658261193086

$ws.Range("C15").Value = 0.3986977119751156
$ws.Range("E15").Value = 1.500453746466346

$ws.Range("C16").Value = 1.97557360987699
$ws.Range("E16").Value = 0.05919149746531627

$ws.Range("C17").Value = -2.311523918755531
$ws.Range("E17").Value = -0.1671238976421296

$ws.Range("C18").Value = -0.8195740704359578
$ws.Range("E18").Value = -0.4351858173977874

$ws.Range("C19").Value = 0.4857399523052974
$ws.Range("E19").Value = -0.3387208997876479
